$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 0.9973287529977295
$ws.Range("E4").Value = 0.01940319488623092

$ws.Range("B5").Value = 68.62320671302086
$ws.Range("E5").Value = 1.241353221023322

$ws.Range("B6").Value = 0.4695827
$ws.Range("C6").Value = 0.5304173
$ws.Range("E6").Value = 0.9981564000000001
$ws.Range("F6").Value = 0.0018436

$ws.Range("B7").Value = 26.5027078
$ws.Range("C7").Value = 26.9435108
$ws.Range("E7").Value = 32.5117888
$ws.Range("F7").Value = 16.068455
